# Bài 2 - Phần 3: Vẽ đường tròn
# Updates the cached "datetimeFigureOut" placeholder text on the slide
# master + every slide layout (2/24/2023 -> 3/2/2023), and renames
# "vòng tròn" to "đường tròn" on slides 6 and 7.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Refresh the cached date/time placeholder text everywhere it lives:
#    the slide master and all of its custom (slide) layouts.
# ---------------------------------------------------------------------
$newDate = "3/2/2023"

$masterShapes = $p.SlideMaster.Shapes
for ($i = 1; $i -le $masterShapes.Count; $i++) {
    $sh = $masterShapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = ""
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

$design = $p.Designs.Item(1)
$layouts = $design.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    $layoutShapes = $layout.Shapes
    for ($i = 1; $i -le $layoutShapes.Count; $i++) {
        $sh = $layoutShapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = ""
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# ---------------------------------------------------------------------
# 2) Slide 6 title: "vòng tròn" -> "đường tròn"
# ---------------------------------------------------------------------
$slide6 = $p.Slides.Item(6)
$title6 = $slide6.Shapes.Item(1).TextFrame.TextRange
$title6.Text = ""
$title6.Text = "5. Vẽ hình đường tròn (Cách 1)"

# ---------------------------------------------------------------------
# 3) Slide 7 title: "vòng tròn" -> "đường tròn"
# ---------------------------------------------------------------------
$slide7 = $p.Slides.Item(7)
$title7 = $slide7.Shapes.Item(1).TextFrame.TextRange
$title7.Text = ""
$title7.Text = "6. Vẽ hình đường tròn (Cách 2)"

# ---------------------------------------------------------------------
# 4) Slide 7 body: second paragraph "vòng tròn" -> "đường tròn"
#    (first paragraph left untouched)
# ---------------------------------------------------------------------
$body7 = $slide7.Shapes.Item(2).TextFrame.TextRange
$body7.Text = ""
$body7.Text = "Cách này là cách vẽ thủ công, dễ hình dung`rTất cả các điểm ở trên đường tròn đều có khoảng cách bằng nhau so với tâm hình tròn"
